# Auto-generated Excel COM-interop script.
# Applies the scheduled market-data refresh to the Pandaemonium Profits workbook:
# updates Universalis price snapshots (currentAveragePrice*) and the dependent
# Leve cost/profit columns for the affected rows across all Grand Company sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1407.675
$ws.Range("I15").Value = 1407.675
$ws.Range("K15").Value = 4223.025
$ws.Range("M15").Value = -4054.025

$ws.Range("H74").Value = 4772.778
$ws.Range("I74").Value = 4173.8184
$ws.Range("J74").Value = 5714
$ws.Range("K74").Value = 4173.8184
$ws.Range("L74").Value = 5714
$ws.Range("M74").Value = -3237.8184
$ws.Range("N74").Value = -7586

$ws.Range("H77").Value = 4772.778
$ws.Range("I77").Value = 4173.8184
$ws.Range("J77").Value = 5714
$ws.Range("K77").Value = 20869.092
$ws.Range("L77").Value = 28570
$ws.Range("M77").Value = -16189.092
$ws.Range("N77").Value = -37930

$ws.Range("H113").Value = 3124.8572
$ws.Range("I113").Value = 1940
$ws.Range("J113").Value = 3783.111
$ws.Range("K113").Value = 1940
$ws.Range("L113").Value = 3783.111
$ws.Range("M113").Value = 1314
$ws.Range("N113").Value = -10291.111

$ws.Range("H137").Value = 2138.2285
$ws.Range("I137").Value = 2225.3447
$ws.Range("J137").Value = 1717.1666
$ws.Range("K137").Value = 6676.034100000001
$ws.Range("L137").Value = 5151.4998
$ws.Range("M137").Value = -4126.034100000001
$ws.Range("N137").Value = -10251.4998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 2277
$ws.Range("I25").Value = 2277
$ws.Range("K25").Value = 2277
$ws.Range("M25").Value = -1875

$ws.Range("H43").Value = 14000
$ws.Range("I43").Value = 9000
$ws.Range("J43").Value = 19000
$ws.Range("K43").Value = 9000
$ws.Range("L43").Value = 19000
$ws.Range("M43").Value = -8687
$ws.Range("N43").Value = -19626

$ws.Range("H101").Value = 40000
$ws.Range("J101").Value = 40000
$ws.Range("L101").Value = 40000
$ws.Range("N101").Value = -46490

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 42798.64
$ws.Range("I134").Value = 2629.0527
$ws.Range("K134").Value = 7887.158100000001
$ws.Range("M134").Value = -5352.158100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5523.1113
$ws.Range("I31").Value = 5161.115
$ws.Range("J31").Value = 6464.3
$ws.Range("K31").Value = 5161.115
$ws.Range("L31").Value = 6464.3
$ws.Range("M31").Value = -4866.115
$ws.Range("N31").Value = -7054.3

$ws.Range("H34").Value = 5523.1113
$ws.Range("I34").Value = 5161.115
$ws.Range("J34").Value = 6464.3
$ws.Range("K34").Value = 5161.115
$ws.Range("L34").Value = 6464.3
$ws.Range("M34").Value = -4959.115
$ws.Range("N34").Value = -6868.3

$ws.Range("H58").Value = 2938652
$ws.Range("I58").Value = 5684426
$ws.Range("K58").Value = 5684426
$ws.Range("M58").Value = -5684223

$ws.Range("H62").Value = 2998.3333
$ws.Range("I62").Value = 2998.3333
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2998.3333
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2374.3333
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 2998.3333
$ws.Range("I65").Value = 2998.3333
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 14991.6665
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -11871.6665
$ws.Range("N65").ClearContents()

$ws.Range("H122").Value = 12162.467
$ws.Range("I122").Value = 5468.5835
$ws.Range("K122").Value = 16405.7505
$ws.Range("M122").Value = -13955.7505

$ws.Range("H132").Value = 13456.818
$ws.Range("I132").Value = 15114.111
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 45342.333
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -42812.333
$ws.Range("N132").Value = -23057

$ws.Range("H134").Value = 3710.72
$ws.Range("I134").Value = 3120.8
$ws.Range("K134").Value = 9362.400000000001
$ws.Range("M134").Value = -6827.400000000001

$ws.Range("H136").Value = 2938652
$ws.Range("I136").Value = 5684426
$ws.Range("K136").Value = 17053278
$ws.Range("M136").Value = -17050728

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 5312.222
$ws.Range("J39").Value = 5312.222
$ws.Range("L39").Value = 15936.666
$ws.Range("N39").Value = -16524.666

$ws.Range("H81").Value = 4441.8667
$ws.Range("J81").Value = 4473.4287
$ws.Range("L81").Value = 13420.2861
$ws.Range("N81").Value = -15666.2861

$ws.Range("H84").Value = 4441.8667
$ws.Range("J84").Value = 4473.4287
$ws.Range("L84").Value = 40260.85830000001
$ws.Range("N84").Value = -51492.85830000001

$ws.Range("H134").Value = 2854.4482

$ws.Range("H137").Value = 104746.6
$ws.Range("I137").Value = 1233.3334
$ws.Range("J137").Value = 260016.5
$ws.Range("K137").Value = 3700.0002
$ws.Range("L137").Value = 780049.5
$ws.Range("M137").Value = 1399.9998
$ws.Range("N137").Value = -790249.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 20714.285
$ws.Range("J57").Value = 20714.285
$ws.Range("L57").Value = 20714.285
$ws.Range("N57").Value = -22354.285

$ws.Range("H126").Value = 1572.909
$ws.Range("J126").Value = 1521.3334
$ws.Range("L126").Value = 4564.0002
$ws.Range("N126").Value = -9504.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 860.6
$ws.Range("I16").Value = 913
$ws.Range("J16").Value = 782
$ws.Range("K16").Value = 913
$ws.Range("L16").Value = 782
$ws.Range("M16").Value = -743
$ws.Range("N16").Value = -1122

$ws.Range("H132").Value = 5642.8335
$ws.Range("I132").Value = 6202.2
$ws.Range("J132").Value = 5243.2856
$ws.Range("K132").Value = 18606.6
$ws.Range("L132").Value = 15729.8568
$ws.Range("M132").Value = -16076.6
$ws.Range("N132").Value = -20789.8568

$ws.Range("H136").Value = 5216.8613
$ws.Range("I136").Value = 2547.4707
$ws.Range("J136").Value = 7605.263
$ws.Range("K136").Value = 7642.4121
$ws.Range("L136").Value = 22815.789
$ws.Range("M136").Value = -5092.4121

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4330.0586
$ws.Range("I132").Value = 3829.1428
$ws.Range("J132").Value = 6667.6665
$ws.Range("K132").Value = 11487.4284
$ws.Range("L132").Value = 20002.9995
$ws.Range("M132").Value = -8957.428400000001
$ws.Range("N132").Value = -25062.9995

$ws.Range("H136").Value = 5666.38
$ws.Range("I136").Value = 2929.3333
$ws.Range("K136").Value = 8787.999899999999
$ws.Range("M136").Value = -6237.999899999999
